# Apply updated Random Forest imputation results to Sheet1.
# (Commit: "Update Name of Algo")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new value, taken from the canonical OOXML diff.
$updates = @{
    "C3"  = -12.3861
    "C4"  = -14.2259
    "D6"  = -7.590799999999997
    "C7"  = -12.01219999999999
    "D7"  = -7.2578
    "C8"  = -12.14249999999999
    "D8"  = -7.686000000000004
    "B11" = 6.041200000000001
    "B12" = 5.345499999999999
    "C12" = -14.67390000000002
    "C14" = -12.15199999999999
    "B15" = 5.9954
    "D19" = -8.549499999999988
    "D21" = -7.703199999999999
    "C22" = -11.10149999999999
    "D24" = -7.299399999999996
    "D25" = -7.726499999999995
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
